$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 729.29034
$ws.Range("I2").Value = 688.6799999999999
$ws.Range("J2").Value = 898.5
$ws.Range("K2").Value = 688.6799999999999
$ws.Range("L2").Value = 898.5
$ws.Range("M2").Value = -575.6799999999999
$ws.Range("N2").Value = -1124.5
$ws.Range("H32").Value = 9962.057000000001
$ws.Range("I32").Value = 6080.1606
$ws.Range("K32").Value = 6080.1606
$ws.Range("M32").Value = -5793.1606
$ws.Range("H61").Value = 971.83636
$ws.Range("I61").Value = 788.76
$ws.Range("J61").Value = 2802.6
$ws.Range("K61").Value = 788.76
$ws.Range("L61").Value = 2802.6
$ws.Range("M61").Value = -576.76
$ws.Range("N61").Value = -3226.6
$ws.Range("H97").Value = 1114.1765
$ws.Range("I97").Value = 1115.7
$ws.Range("J97").Value = 1112
$ws.Range("K97").Value = 1115.7
$ws.Range("L97").Value = 1112
$ws.Range("M97").Value = -619.7
$ws.Range("N97").Value = -2104
$ws.Range("H116").Value = 729.29034
$ws.Range("I116").Value = 688.6799999999999
$ws.Range("J116").Value = 898.5
$ws.Range("K116").Value = 688.6799999999999
$ws.Range("L116").Value = 898.5
$ws.Range("M116").Value = 1605.32
$ws.Range("N116").Value = -5486.5
$ws.Range("H132").Value = 1925.1041
$ws.Range("I132").Value = 901.4545000000001
$ws.Range("K132").Value = 2704.3635
$ws.Range("M132").Value = -174.3635000000004
$ws.Range("H136").Value = 971.83636
$ws.Range("I136").Value = 788.76
$ws.Range("J136").Value = 2802.6
$ws.Range("K136").Value = 2366.28
$ws.Range("L136").Value = 8407.799999999999
$ws.Range("M136").Value = 183.7200000000003
$ws.Range("N136").Value = -13507.8
$ws.Range("H139").Value = 43330
$ws.Range("J139").Value = 43330
$ws.Range("L139").Value = 43330
$ws.Range("N139").Value = -53610

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 729.29034
$ws.Range("I3").Value = 688.6799999999999
$ws.Range("J3").Value = 898.5
$ws.Range("K3").Value = 688.6799999999999
$ws.Range("L3").Value = 898.5
$ws.Range("M3").Value = -574.6799999999999
$ws.Range("N3").Value = -1126.5
$ws.Range("H99").Value = 3325.8845
$ws.Range("I99").Value = 1866.25
$ws.Range("J99").Value = 4577
$ws.Range("K99").Value = 1866.25
$ws.Range("L99").Value = 4577
$ws.Range("M99").Value = -368.25
$ws.Range("N99").Value = -7573
$ws.Range("H134").Value = 2680.228
$ws.Range("I134").Value = 1584.3617
$ws.Range("J134").Value = 7830.8
$ws.Range("K134").Value = 4753.0851
$ws.Range("L134").Value = 23492.4
$ws.Range("M134").Value = -2218.0851
$ws.Range("N134").Value = -28562.4
$ws.Range("H138").Value = 41248.78
$ws.Range("J138").Value = 41248.78
$ws.Range("L138").Value = 41248.78
$ws.Range("N138").Value = -51528.78

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2493.0425
$ws.Range("I31").Value = 970.9706
$ws.Range("K31").Value = 970.9706
$ws.Range("M31").Value = -675.9706
$ws.Range("H34").Value = 2493.0425
$ws.Range("I34").Value = 970.9706
$ws.Range("K34").Value = 970.9706
$ws.Range("M34").Value = -768.9706
$ws.Range("H58").Value = 2008.5735
$ws.Range("I58").Value = 1687.7869
$ws.Range("J58").Value = 4804
$ws.Range("K58").Value = 1687.7869
$ws.Range("L58").Value = 4804
$ws.Range("M58").Value = -1484.7869
$ws.Range("N58").Value = -5210
$ws.Range("H132").Value = 2731.8
$ws.Range("I132").Value = 1772.4634
$ws.Range("K132").Value = 5317.3902
$ws.Range("M132").Value = -2787.3902
$ws.Range("H134").Value = 3355.7
$ws.Range("I134").Value = 3382.5527
$ws.Range("J134").Value = 3270.6667
$ws.Range("K134").Value = 10147.6581
$ws.Range("L134").Value = 9812.000100000001
$ws.Range("M134").Value = -7612.658100000001
$ws.Range("N134").Value = -14882.0001
$ws.Range("H136").Value = 2008.5735
$ws.Range("I136").Value = 1687.7869
$ws.Range("J136").Value = 4804
$ws.Range("K136").Value = 5063.3607
$ws.Range("L136").Value = 14412
$ws.Range("M136").Value = -2513.3607
$ws.Range("N136").Value = -19512
$ws.Range("H138").Value = 47447.6
$ws.Range("J138").Value = 47447.6
$ws.Range("L138").Value = 47447.6
$ws.Range("N138").Value = -57727.6
$ws.Range("H140").Value = 92824.164
$ws.Range("J140").Value = 92824.164
$ws.Range("L140").Value = 92824.164
$ws.Range("N140").Value = -103184.164
$ws.Range("H141").Value = 15381.579
$ws.Range("J141").Value = 15381.579
$ws.Range("L141").Value = 15381.579
$ws.Range("N141").Value = -25741.579

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3336656.2
$ws.Range("I87").Value = 3336656.2
$ws.Range("K87").Value = 10009968.6
$ws.Range("M87").Value = -10008720.6
$ws.Range("H90").Value = 3336656.2
$ws.Range("I90").Value = 3336656.2
$ws.Range("K90").Value = 30029905.8
$ws.Range("M90").Value = -30023665.8
$ws.Range("H131").Value = 13514527
$ws.Range("I131").Value = 71430264
$ws.Range("J131").Value = 856.2
$ws.Range("K131").Value = 214290792
$ws.Range("L131").Value = 2568.6
$ws.Range("M131").Value = -214285752
$ws.Range("N131").Value = -12648.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2377.2983
$ws.Range("I132").Value = 1646.5834
$ws.Range("J132").Value = 3629.9524
$ws.Range("K132").Value = 4939.7502
$ws.Range("L132").Value = 10889.8572
$ws.Range("M132").Value = -2409.7502
$ws.Range("N132").Value = -15949.8572
$ws.Range("H140").Value = 42498.094
$ws.Range("J140").Value = 42498.094
$ws.Range("L140").Value = 42498.094
$ws.Range("N140").Value = -52858.094

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3899.6128
$ws.Range("I122").Value = 3321
$ws.Range("K122").Value = 9963
$ws.Range("M122").Value = -7513
$ws.Range("H139").Value = 47565
$ws.Range("J139").Value = 47565
$ws.Range("L139").Value = 47565
$ws.Range("N139").Value = -57845
$ws.Range("H140").Value = 69972.75
$ws.Range("J140").Value = 69972.75
$ws.Range("L140").Value = 69972.75
$ws.Range("N140").Value = -80332.75
$ws.Range("H141").Value = 41478.75
$ws.Range("J141").Value = 41977.824
$ws.Range("L141").Value = 41977.824
$ws.Range("N141").Value = -52337.824

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7250422.5
$ws.Range("I132").Value = 6628.8237
$ws.Range("J132").Value = 11496784
$ws.Range("K132").Value = 19886.4711
$ws.Range("L132").Value = 34490352
$ws.Range("M132").Value = -17356.4711
$ws.Range("N132").Value = -34495412
$ws.Range("H136").Value = 2697.5854
$ws.Range("I136").Value = 830.8461
$ws.Range("J136").Value = 5933.2666
$ws.Range("K136").Value = 2492.5383
$ws.Range("L136").Value = 17799.7998
$ws.Range("M136").Value = 57.46169999999984
$ws.Range("N136").Value = -22899.7998
$ws.Range("H139").Value = 40359.285
$ws.Range("J139").Value = 40746.54
$ws.Range("L139").Value = 40746.54
$ws.Range("N139").Value = -51026.54
$ws.Range("H140").Value = 38296.668
$ws.Range("J140").Value = 38296.668
$ws.Range("L140").Value = 38296.668
$ws.Range("N140").Value = -48656.668
$ws.Range("H141").Value = 44017.918
$ws.Range("J141").Value = 44017.918
$ws.Range("L141").Value = 44017.918
$ws.Range("N141").Value = -54377.918
